# Updated cryptos list with refreshed Price / Volume(1h) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (preventing Excel from
# auto-converting numeric-looking strings like "1.170" or "22.407.87"
# into numbers), while keeping the cell style/format unchanged.
function Set-TextValue($address, $value) {
    $c = $ws.Range($address)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "22.407.87"
Set-TextValue "E2" "  -4.57%  "
Set-TextValue "D3" "1.570.54"
Set-TextValue "E3" "  -4.73%  "
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "E5" "  -0.05%  "
Set-TextValue "D6" "291.89"
Set-TextValue "E6" "  -2.55%  "
Set-TextValue "D7" "0.3694"
Set-TextValue "E7" "  -2.58%  "
Set-TextValue "D8" "49.71"
Set-TextValue "E8" "  -1.07%  "
Set-TextValue "D9" "0.3369"
Set-TextValue "E9" "  -5.50%  "
Set-TextValue "D10" "1.170"
Set-TextValue "E10" "  -4.12%  "
Set-TextValue "D11" "0.07571"
Set-TextValue "E11" "  -6.47%  "
Set-TextValue "E12" "  -0.11%  "
Set-TextValue "D13" "21.12"
Set-TextValue "E13" "  -4.05%  "
Set-TextValue "D14" "6.059"
Set-TextValue "E14" "  -5.32%  "
Set-TextValue "D15" "6.858"
Set-TextValue "E15" "  -7.12%  "
Set-TextValue "D16" "0.00001144"
Set-TextValue "E16" "  -4.37%  "
Set-TextValue "D17" "1.570.22"
Set-TextValue "E17" "  -5.28%  "
Set-TextValue "D18" "89.31"
Set-TextValue "E18" "  -8.15%  "
Set-TextValue "D19" "0.06703"
Set-TextValue "E19" "  -3.73%  "
Set-TextValue "D20" "1.000"
Set-TextValue "D21" "6.252"
Set-TextValue "E21" "  -7.31%  "
Set-TextValue "D22" "16.37"
Set-TextValue "E22" "  -5.21%  "
Set-TextValue "D23" "11.95"
Set-TextValue "E23" "  -3.66%  "
Set-TextValue "D24" "22.415.59"
Set-TextValue "E24" "  -4.60%  "
Set-TextValue "D25" "2.404"
Set-TextValue "E25" "  -4.31%  "
Set-TextValue "D26" "2.962"
Set-TextValue "E26" "  +1.86%  "
Set-TextValue "D27" "19.85"
Set-TextValue "E27" "  -5.05%  "
Set-TextValue "D28" "146.26"
Set-TextValue "E28" "  -4.61%  "
Set-TextValue "D29" "4.925"
Set-TextValue "E29" "  -5.61%  "
Set-TextValue "D30" "125.31"
Set-TextValue "E30" "  -5.71%  "
Set-TextValue "D31" "1.750.60"
Set-TextValue "E31" "  -4.63%  "
Set-TextValue "D32" "6.286"
Set-TextValue "E32" "  -9.35%  "
Set-TextValue "D33" "1.974"
Set-TextValue "E33" "  -6.72%  "
Set-TextValue "D34" "0.9864"
Set-TextValue "E34" "  -3.18%  "
Set-TextValue "D35" "10.37"
Set-TextValue "E35" "  -12.86%  "
Set-TextValue "D36" "0.08432"
Set-TextValue "E36" "  -3.49%  "
Set-TextValue "D37" "0.02541"
Set-TextValue "E37" "  -6.88%  "
Set-TextValue "D38" "0.2306"
Set-TextValue "E38" "  -5.26%  "
Set-TextValue "D39" "0.06521"
Set-TextValue "E39" "  -4.32%  "
Set-TextValue "D40" "5.505"
Set-TextValue "E40" "  -7.47%  "
Set-TextValue "D41" "11.80"
Set-TextValue "E41" "  -10.36%  "
Set-TextValue "D42" "1.246"
Set-TextValue "E42" "  -5.28%  "
Set-TextValue "D43" "0.6403"
Set-TextValue "E43" "  -7.16%  "
Set-TextValue "D44" "14.54"
Set-TextValue "E44" "  -6.29%  "
Set-TextValue "D45" "0.9996"
Set-TextValue "E45" "  -0.12%  "
Set-TextValue "D46" "0.6028"
Set-TextValue "E46" "  -6.02%  "
Set-TextValue "D47" "3.776"
Set-TextValue "E47" "  -3.73%  "
Set-TextValue "D48" "2.113"
Set-TextValue "E48" "  -6.62%  "
Set-TextValue "D49" "121.83"
Set-TextValue "E49" "  -4.47%  "
Set-TextValue "D50" "0.07265"
Set-TextValue "E50" "  -6.20%  "
Set-TextValue "D51" "1.189"
Set-TextValue "E51" "  +0.54%  "
